$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: estoque_atualizado (G4) and desvio_padrao (I4)
$ws.Range("G4").Value = -21
$ws.Range("I4").Value = 0.22

# Row 6: estoque_atualizado (G6)
$ws.Range("G6").Value = -88

# Row 7: estoque_atualizado (G7)
$ws.Range("G7").Value = 8

# Row 9: estoque_atualizado (G9), media_vendas (H9), desvio_padrao (I9)
$ws.Range("G9").Value = -50
$ws.Range("H9").Value = 1.03
$ws.Range("I9").Value = 0.18
